$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "pixel_size_mm" in I1, formatted bold like the other headers
$ws.Range("I1").Value = "pixel_size_mm"
$ws.Range("I1").Font.Bold = $true

# Add the new pixel size value in I2 (DVF errors multiplied by pixel size)
$ws.Range("I2").Value = 1.8180000000000001

# Update the active selection as in the saved workbook
$ws.Range("I5").Select()

$wb.Save()
